$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 73

# Column A holds a date-formatted string ("2020-08-11"). A direct literal
# Value assignment gets auto-detected by Excel as a date and rewritten as a
# serial number with a date number format (adding an unwanted style). Using
# a formula that evaluates to the text, then converting the formula result
# to a plain value via copy / PasteSpecial(values), keeps it as plain text
# (stored as a shared string) without touching the style table.
$ws.Cells.Item($row, 1).Formula = '="2020-08-11"'
$ws.Cells.Item($row, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 2).Value = 492522
$ws.Cells.Item($row, 3).Value = 538333
$ws.Cells.Item($row, 4).Value = 81259
$ws.Cells.Item($row, 5).Value = 53929
$ws.Cells.Item($row, 6).Value = 26.51
